$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row - new columns I and J
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the style from H1 (existing header cell) to I1 and J1 so they match
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1:J1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$ws.Application.CutCopyMode = $false

# Data values for columns I (I0) and J (IF)
$values = @{
    2  = @(8, 8)
    3  = @(6, 7)
    4  = @(6, 7)
    5  = @(5, 6)
    6  = @(9, 9)
    7  = @(5, 6)
    8  = @(7, 7)
    9  = @(3, 3)
    10 = @(7, 7)
}

foreach ($row in $values.Keys) {
    $pair = $values[$row]
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
}
